$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated coin price (D) / 1h volume change (E) figures refreshed by the scraper run.
$ws.Range("D2").Value = "'308.50"
$ws.Range("E2").Value = "'-0.05%"
$ws.Range("D2:E2").Style = "Normal"

$ws.Range("D3").Value = "'40.73"
$ws.Range("E3").Value = "'1.67%"
$ws.Range("D3:E3").Style = "Normal"

$ws.Range("D4").Value = "'5.117"
$ws.Range("E4").Value = "'-0.05%"
$ws.Range("D4:E4").Style = "Normal"

$ws.Range("D5").Value = "'0.07636"
$ws.Range("E5").Value = "'-1.41%"
$ws.Range("D5:E5").Style = "Normal"

$ws.Range("E6").Value = "'0.53%"
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'0.9032"
$ws.Range("E7").Value = "'2.31%"
$ws.Range("D7:E7").Style = "Normal"

$ws.Range("D9").Value = "'0.1116"
$ws.Range("E9").Value = "'10.34%"
$ws.Range("D9:E9").Style = "Normal"

$ws.Range("D10").Value = "'0.1784"
$ws.Range("E10").Value = "'2.09%"
$ws.Range("D10:E10").Style = "Normal"

$ws.Range("D11").Value = "'0.09174"
$ws.Range("E11").Value = "'1.50%"
$ws.Range("D11:E11").Style = "Normal"

$ws.Range("D12").Value = "'0.04209"
$ws.Range("E12").Value = "'-5.15%"
$ws.Range("D12:E12").Style = "Normal"

$ws.Range("D13").Value = "'0.1052"
$ws.Range("E13").Value = "'-0.39%"
$ws.Range("D13:E13").Style = "Normal"

$ws.Range("D14").Value = "'0.001251"
$ws.Range("E14").Value = "'-0.88%"
$ws.Range("D14:E14").Style = "Normal"

$ws.Range("D15").Value = "'0.005876"
$ws.Range("E15").Value = "'0.56%"
$ws.Range("D15:E15").Style = "Normal"

$ws.Range("D16").Value = "'3.351"
$ws.Range("E16").Value = "'-0.16%"
$ws.Range("D16:E16").Style = "Normal"

$ws.Range("D17").Value = "'4.242"
$ws.Range("E17").Value = "'-0.36%"
$ws.Range("D17:E17").Style = "Normal"

$ws.Range("E18").Value = "'0.49%"
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'6.619"
$ws.Range("E19").Value = "'-6.55%"
$ws.Range("D19:E19").Style = "Normal"

$ws.Range("D20").Value = "'0.1364"
$ws.Range("E20").Value = "'1.77%"
$ws.Range("D20:E20").Style = "Normal"

$ws.Range("D21").Value = "'0.2790"
$ws.Range("E21").Value = "'-0.12%"
$ws.Range("D21:E21").Style = "Normal"

$ws.Range("D22").Value = "'0.04074"
$ws.Range("E22").Value = "'-2.37%"
$ws.Range("D22:E22").Style = "Normal"

$ws.Range("E23").Value = "'2.06%"
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'0.004122"
$ws.Range("E24").Value = "'0.71%"
$ws.Range("D24:E24").Style = "Normal"

$ws.Range("D25").Value = "'0.0001301"
$ws.Range("E25").Value = "'-0.22%"
$ws.Range("D25:E25").Style = "Normal"

$ws.Range("D26").Value = "'0.0003746"
$ws.Range("D26").Style = "Normal"

$ws.Range("D38").Value = "'0.02417"
$ws.Range("E38").Value = "'2.56%"
$ws.Range("D38:E38").Style = "Normal"

$ws.Range("D39").Value = "'0.05177"
$ws.Range("E39").Value = "'-0.97%"
$ws.Range("D39:E39").Style = "Normal"

$ws.Range("D40").Value = "'0.007784"
$ws.Range("E40").Value = "'-2.00%"
$ws.Range("D40:E40").Style = "Normal"

$ws.Range("D41").Value = "'0.1302"
$ws.Range("E41").Value = "'-1.95%"
$ws.Range("D41:E41").Style = "Normal"

$ws.Range("D42").Value = "'0.007047"
$ws.Range("E42").Value = "'10.90%"
$ws.Range("D42:E42").Style = "Normal"

$ws.Range("D43").Value = "'0.001951"
$ws.Range("E43").Value = "'-0.57%"
$ws.Range("D43:E43").Style = "Normal"

$ws.Range("D44").Value = "'0.007973"
$ws.Range("E44").Value = "'-9.17%"
$ws.Range("D44:E44").Style = "Normal"

$ws.Range("E45").Value = "'-7.71%"
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'0.00006996"
$ws.Range("E46").Value = "'6.50%"
$ws.Range("D46:E46").Style = "Normal"

$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("E47").Value = "'-0.22%"
$ws.Range("D47:E47").Style = "Normal"

$ws.Range("D48").Value = "'0.03147"
$ws.Range("E48").Value = "'662.94%"
$ws.Range("D48:E48").Style = "Normal"

$ws.Range("D50").Value = "'0.00002101"
$ws.Range("E50").Value = "'-0.22%"
$ws.Range("D50:E50").Style = "Normal"

$ws.Range("D51").Value = "'0.0002001"
$ws.Range("E51").Value = "'-0.22%"
$ws.Range("D51:E51").Style = "Normal"
